$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- row 75: extend the array formula with a third IF argument -----------
$ws.Range("C75").FormulaArray = "=SUM(IF(A67:B68>2,A67:B68,1))"

# --- new block: "Use array cell" ------------------------------------------
# A97:B98 is a 2x2 array formula (SQRT of a 2x2 literal array)
$ws.Range("A97:B98").FormulaArray = "=SQRT({1,2;3,4})"

# C99 references a single cell (B98) that is part of the array above
$ws.Range("C99").Formula = "=B98+1"
$ws.Range("I99").Value = 3
$ws.Range("N99").Value = "Use array cell"

# --- new block: "ref array in array" --------------------------------------
# B102:C103 is a 2x2 array formula referencing the A97:B98 array range
$ws.Range("B102:C103").FormulaArray = "=A97:B98+1"
$ws.Range("I102").Value = 2
$ws.Range("J102").Value = 2.4142135623730949
$ws.Range("N102").Value = "ref array in array"
$ws.Range("I103").Value = 2.7320508075688772
$ws.Range("J103").Value = 3

# --- move the trailing "END" marker row from 101 to 110 -------------------
$ws.Range("H101").ClearContents()
$ws.Range("H110").Value = "END"
$ws.Range("I110").Value = "END"

# --- column C gets an explicit width ---------------------------------------
$ws.Columns.Item(3).ColumnWidth = 12.1

# --- view state: scroll position + selection -------------------------------
$excel.ActiveWindow.ScrollRow = 88
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("J103").Select()
